$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04193501389126
$ws.Range("D2").Value = 1.048556341196559
$ws.Range("E2").Value = 1.055020396193695
$ws.Range("F2").Value = 1.061389353261513
$ws.Range("I2").Value = 1.041418978249233
$ws.Range("J2").Value = 1.047013544554417
$ws.Range("K2").Value = 1.051315856538892
$ws.Range("L2").Value = 1.057762003201966
$ws.Range("M2").Value = 1.064113546403013
$ws.Range("N2").Value = 1.019566382530981
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042748953802618
$ws.Range("D3").Value = 1.049185031758192
$ws.Range("E3").Value = 1.055794751089229
$ws.Range("F3").Value = 1.062185263928952
$ws.Range("I3").Value = 1.041590927976909
$ws.Range("J3").Value = 1.047474391361324
$ws.Range("K3").Value = 1.051757090794496
$ws.Range("L3").Value = 1.058349816303406
$ws.Range("M3").Value = 1.06472411563336
$ws.Range("N3").Value = 1.01972083635055
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043276400249202
$ws.Range("D4").Value = 1.049592509500735
$ws.Range("E4").Value = 1.056296924740789
$ws.Range("F4").Value = 1.062701423078262
$ws.Range("I4").Value = 1.041701374638252
$ws.Range("J4").Value = 1.047772675585868
$ws.Range("K4").Value = 1.05204257923944
$ws.Range("L4").Value = 1.058730621478928
$ws.Range("M4").Value = 1.065119686539399
$ws.Range("N4").Value = 1.019820768565562
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043498321939337
$ws.Range("D5").Value = 1.049763972350542
$ws.Range("E5").Value = 1.056508303253042
$ws.Range("F5").Value = 1.062918690007118
$ws.Range("I5").Value = 1.041747610246543
$ws.Range("J5").Value = 1.04789809339816
$ws.Range("K5").Value = 1.052162592503417
$ws.Range("L5").Value = 1.05889081824318
$ws.Range("M5").Value = 1.065286100556274
$ws.Range("N5").Value = 1.019862777241466
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043535594257663
$ws.Range("D6").Value = 1.049792770998494
$ws.Range("E6").Value = 1.056543810098537
$ws.Range("F6").Value = 1.062955186064683
$ws.Range("I6").Value = 1.041755361883698
$ws.Range("J6").Value = 1.047919152682909
$ws.Range("K6").Value = 1.052182742861951
$ws.Range("L6").Value = 1.058917722168389
$ws.Range("M6").Value = 1.065314048954519
$ws.Range("N6").Value = 1.019869830502751
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043279364860108
$ws.Range("D7").Value = 1.049594799972529
$ws.Range("E7").Value = 1.05629974815326
$ws.Range("F7").Value = 1.062704325137092
$ws.Range("I7").Value = 1.041701993212443
$ws.Range("J7").Value = 1.047774351351362
$ws.Range("K7").Value = 1.052044182888162
$ws.Range("L7").Value = 1.058732761620444
$ws.Range("M7").Value = 1.065121909717973
$ws.Range("N7").Value = 1.019821329899279
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042209927560993
$ws.Range("D8").Value = 1.048768669649701
$ws.Range("E8").Value = 1.055281861112676
$ws.Range("F8").Value = 1.06165809526113
$ws.Range("I8").Value = 1.041477258040727
$ws.Range("J8").Value = 1.047169271370045
$ws.Range("K8").Value = 1.051464976954459
$ws.Range("L8").Value = 1.057960563058683
$ws.Range("M8").Value = 1.064319788267362
$ws.Range("N8").Value = 1.019618582639347
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040331442906673
$ws.Range("D9").Value = 1.047318159685612
$ws.Range("E9").Value = 1.053496838457927
$ws.Range("F9").Value = 1.05982341844811
$ws.Range("I9").Value = 1.041075026878425
$ws.Range("J9").Value = 1.046103757911817
$ws.Range("K9").Value = 1.050444252894094
$ws.Range("L9").Value = 1.056603376177305
$ws.Range("M9").Value = 1.062910190685089
$ws.Range("N9").Value = 1.019261263482461
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039083263616651
$ws.Range("D10").Value = 1.046354787245632
$ws.Range("E10").Value = 1.052312744290016
$ws.Range("F10").Value = 1.058606419633478
$ws.Range("I10").Value = 1.04080273589358
$ws.Range("J10").Value = 1.045393976959121
$ws.Range("K10").Value = 1.04976379280711
$ws.Range("L10").Value = 1.055701054216651
$ws.Range("M10").Value = 1.061973146023298
$ws.Range("N10").Value = 1.019023045055955
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038543794186932
$ws.Range("D11").Value = 1.04593852229004
$ws.Range("E11").Value = 1.051801448759494
$ws.Range("F11").Value = 1.058080923755773
$ws.Range("I11").Value = 1.04068385908977
$ws.Range("J11").Value = 1.045086784216197
$ws.Range("K11").Value = 1.049469168935775
$ws.Range("L11").Value = 1.055310945470597
$ws.Range("M11").Value = 1.061568054565845
$ws.Range("N11").Value = 1.018919898792719
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.038343563395898
$ws.Range("D12").Value = 1.045784037364165
$ws.Range("E12").Value = 1.051611747048535
$ws.Range("F12").Value = 1.057885954612519
$ws.Range("I12").Value = 1.040639557517668
$ws.Range("J12").Value = 1.044972702636244
$ws.Range("K12").Value = 1.049359736907493
$ws.Range("L12").Value = 1.055166133949298
$ws.Range("M12").Value = 1.061417685756084
$ws.Range("N12").Value = 1.018881586770637
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.038386506663974
$ws.Range("D13").Value = 1.045817168791398
$ws.Range("E13").Value = 1.051652428891378
$ws.Range("F13").Value = 1.057927766010283
$ws.Range("I13").Value = 1.040649066924418
$ws.Range("J13").Value = 1.044997172437558
$ws.Range("K13").Value = 1.049383210227401
$ws.Range("L13").Value = 1.055197192311001
$ws.Range("M13").Value = 1.061449935813674
$ws.Range("N13").Value = 1.01888980477182
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038527239932131
$ws.Range("D14").Value = 1.045925749764615
$ws.Range("E14").Value = 1.051785763533949
$ws.Range("F14").Value = 1.058064802975223
$ws.Range("I14").Value = 1.040680200070471
$ws.Range("J14").Value = 1.045077353710173
$ws.Range("K14").Value = 1.049460123152325
$ws.Range("L14").Value = 1.055298973413645
$ws.Range("M14").Value = 1.06155562297891
$ws.Range("N14").Value = 1.018916731883792
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038613970564092
$ws.Range("D15").Value = 1.045992667955186
$ws.Range("E15").Value = 1.051867944121613
$ws.Range("F15").Value = 1.058149265642945
$ws.Range("I15").Value = 1.040699362980295
$ws.Range("J15").Value = 1.04512675919159
$ws.Range("K15").Value = 1.049507512368618
$ws.Range("L15").Value = 1.055361696383269
$ws.Range("M15").Value = 1.061620753659412
$ws.Range("N15").Value = 1.018933322726181
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.039119087663693
$ws.Range("D16").Value = 1.04638243212223
$ws.Range("E16").Value = 1.052346707495423
$ws.Range("F16").Value = 1.058641326276801
$ws.Range("I16").Value = 1.040810604919272
$ws.Range("J16").Value = 1.045414367546928
$ws.Range("K16").Value = 1.049783346561892
$ws.Range("L16").Value = 1.055726957290923
$ws.Range("M16").Value = 1.062000044562451
$ws.Range("N16").Value = 1.019029890659072
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039436203440683
$ws.Range("D17").Value = 1.046627158313298
$ws.Range("E17").Value = 1.052647406133463
$ws.Range("F17").Value = 1.058950378716339
$ws.Range("I17").Value = 1.04088012407468
$ws.Range("J17").Value = 1.045594817027361
$ws.Range("K17").Value = 1.049956376373991
$ws.Range("L17").Value = 1.055956238346316
$ws.Range("M17").Value = 1.062238140360052
$ws.Range("N17").Value = 1.019090466572048
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039621268110549
$ws.Range("D18").Value = 1.046769987897528
$ws.Range("E18").Value = 1.052822935917236
$ws.Range("F18").Value = 1.059130785718118
$ws.Range("I18").Value = 1.040920579490076
$ws.Range("J18").Value = 1.045700084255937
$ws.Range("K18").Value = 1.050057303554994
$ws.Range("L18").Value = 1.056090032136639
$ws.Range("M18").Value = 1.062377080695551
$ws.Range("N18").Value = 1.019125799818218
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039684386675123
$ws.Range("D19").Value = 1.046818703442135
$ws.Range("E19").Value = 1.052882810270919
$ws.Range("F19").Value = 1.059192323853292
$ws.Range("I19").Value = 1.040934357786024
$ws.Range("J19").Value = 1.045735980009984
$ws.Range("K19").Value = 1.050091717381
$ws.Range("L19").Value = 1.056135662152736
$ws.Range("M19").Value = 1.062424466405624
$ws.Range("N19").Value = 1.019137847584428
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039402169934137
$ws.Range("D20").Value = 1.046600892709508
$ws.Range("E20").Value = 1.052615129792794
$ws.Range("F20").Value = 1.058917205616947
$ws.Range("I20").Value = 1.040872675037632
$ws.Range("J20").Value = 1.045575455037505
$ws.Range("K20").Value = 1.049937811714333
$ws.Range("L20").Value = 1.055931632644409
$ws.Range("M20").Value = 1.062212588396285
$ws.Range("N20").Value = 1.019083967310932
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.038485793276532
$ws.Range("D21").Value = 1.045893771648196
$ws.Range("E21").Value = 1.051746493819428
$ws.Range("F21").Value = 1.058024442819556
$ws.Range("I21").Value = 1.040671036143598
$ws.Range("J21").Value = 1.045053741664511
$ws.Range("K21").Value = 1.049437474080231
$ws.Range("L21").Value = 1.055268998855911
$ws.Range("M21").Value = 1.061524497962803
$ws.Range("N21").Value = 1.018908802486211
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.037910512266877
$ws.Range("D22").Value = 1.045449955206926
$ws.Range("E22").Value = 1.051201599903616
$ws.Range("F22").Value = 1.057464421150328
$ws.Range("I22").Value = 1.040543416703013
$ws.Range("J22").Value = 1.044725856307642
$ws.Range("K22").Value = 1.049122918083698
$ws.Range("L22").Value = 1.054852909186806
$ws.Range("M22").Value = 1.061092448769707
$ws.Range("N22").Value = 1.018798675981952
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.038215395446146
$ws.Range("D23").Value = 1.045685156242927
$ws.Range("E23").Value = 1.051490339097212
$ws.Range("F23").Value = 1.057761175875783
$ws.Range("I23").Value = 1.040611149655529
$ws.Range("J23").Value = 1.044899661134035
$ws.Range("K23").Value = 1.049289667219446
$ws.Range("L23").Value = 1.055073434961688
$ws.Range("M23").Value = 1.061321430596433
$ws.Range("N23").Value = 1.018857055332723
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039417547913164
$ws.Range("D24").Value = 1.046612760743092
$ws.Range("E24").Value = 1.052629713657702
$ws.Range("F24").Value = 1.058932194674357
$ws.Range("I24").Value = 1.040876041227165
$ws.Range("J24").Value = 1.04558420384391
$ws.Range("K24").Value = 1.049946200279492
$ws.Range("L24").Value = 1.055942750722884
$ws.Range("M24").Value = 1.062224134034154
$ws.Range("N24").Value = 1.019086904046286
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040816353863774
$ws.Range("D25").Value = 1.047692518991554
$ws.Range("E25").Value = 1.053957274953789
$ws.Range("F25").Value = 1.060296657066327
$ws.Range("I25").Value = 1.041179745691936
$ws.Range("J25").Value = 1.046379125653846
$ws.Range("K25").Value = 1.050708136590595
$ws.Range("L25").Value = 1.056953813524404
$ws.Range("M25").Value = 1.06327413939829
$ws.Range("N25").Value = 1.01935364221924
